$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (refresh from source), preserving
# each cell as literal text (matching the original inline-string storage)
# by forcing a Text number format before writing the value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.07%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.11%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.023"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.90%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08081"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.98%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.949"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.83%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.141"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.44%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.839"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.89%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9319"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.35%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1253"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-20.30%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1911"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09234"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.98%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03506"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.43%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09931"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.56%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001416"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.73%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006698"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "16.36%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.614"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.25%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.086"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "7.45%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3442"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.02%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.178"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.95%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1304"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.28%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2531"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.41%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04407"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.90%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.98%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004725"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.13%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.74%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003130"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.63%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01968"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.19%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05167"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "8.04%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007603"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.39%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.38%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1371"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.99%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.43%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01070"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "10.22%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006374"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.77%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.17%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.22"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.85%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001601"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.54%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.17%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.17%"
